# Apply the "contributions-summary" edit: fill in real usernames for the
# first three group members (previously placeholder "a1DDDDDDD") and fill
# in the assessment grade/comment cells for Iteration 2 and Iteration 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Iteration 3 & Iteration 2 comments first (matches shared-string / ---
# --- style creation order observed in the target workbook) -------------
$ws.Range("D33").Value = "fixed some style issues"
$ws.Range("D25").Value = "mostly completed voting page."
$ws.Range("D23").Value = "Strong communication and work ethic, continues to help other memebers "
$ws.Range("D22").Value = "Great team work and overall coding"
$ws.Range("D24").Value = "Strong communication and suggestions on path of project "
$ws.Range("D31").Value = "Helped to set the unique style of the whole application."
$ws.Range("D30").Value = "Finalized project."
$ws.Range("D32").Value = "Finalized project."

# Grades (HD) for Iteration 2 and Iteration 3 rows
$ws.Range("C22").Value = "HD"
$ws.Range("C23").Value = "HD"
$ws.Range("C24").Value = "HD"
$ws.Range("C25").Value = "HD"
$ws.Range("C30").Value = "HD"
$ws.Range("C31").Value = "HD"
$ws.Range("C32").Value = "HD"
$ws.Range("C33").Value = "HD"

# --- Iteration 1 usernames: replace placeholders with real usernames ---
$ws.Range("A6").Value = "a1724402 "
$ws.Range("A7").Value = "a1720458"
$ws.Range("A8").Value = "a1725532"
$ws.Range("A9").Value = "a1702065"

# Row 8's username cell picks up a bottom border (mirrors the source
# workbook, which only thickens the border under column A here).
$ws.Range("A8").Borders.Item(9).LineStyle = 1
$ws.Range("A8").Borders.Item(9).Weight = -4138
$ws.Rows.Item(8).RowHeight = 16.15

# Row 24's comment cell (D24) picks up the thicker bottom border instead.
$ws.Range("D24").Borders.Item(9).LineStyle = 1
$ws.Range("D24").Borders.Item(9).Weight = -4138
$ws.Rows.Item(24).RowHeight = 16.15

# Final selection left on A9, matching the saved view state.
$ws.Range("A9").Select()
